# Change network line settings:
# The R + L + (C // G) branch has been changed to (R + L)//C//G branch.
#
# On the "NetworkLine" sheet:
#   - Mutual branches (rows 10-12): G (pu) column (F) goes from "inf" to 0
#   - Self branches (rows 13-18): R (pu) and wL (pu) columns (C, D) go from 0 to "inf"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NetworkLine")

# Mutual branches: remove the shunt G (now a closed/zero conductance)
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 0

# Self branches: open up the series R and L (now "inf")
$ws.Range("C13:D18").Value = "inf"

# Reflect that the user ended up on the NetworkLine sheet with D14 selected
$ws.Activate()
$ws.Range("D14").Select() | Out-Null
